# Generate Report for Handoff
# Updates the localization-status report: flips the "In Translation" status
# to "Ready for handoff" and bumps the associated handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 06:44:59"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 06:44:55"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 06:44:59"

# --- Widen the Status/date columns to fit the new, longer text ----------
# (ColumnWidth is quantized by the host to whole pixels, so we pick the
#  character-width input that lands closest to the target OOXML column
#  width of 17.2159881591797.)
$wsOverview.Columns("E:F").ColumnWidth = 16.3333333333333
$wsZhCn.Columns("C:C").ColumnWidth = 16.3333333333333
$wsDeDe.Columns("C:C").ColumnWidth = 16.3333333333333
